$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 12
$ws.Range("D3").Value = 2
$ws.Range("D4").Value = 2
$ws.Range("D5").Value = 2

$ws.Rows.Item(1).RowHeight = 19.5
